$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# Add a new translation row ("nav5" / "R code" / "Code R") directly below
# the existing "nav4" row, and repurpose the "nav4" row's English/French
# text to describe the new "4. Report" tab (the translation table's ids
# track tab position, so the tab that used to be "nav4" - R code - becomes
# "nav5", and the newly inserted "Report" tab takes over the "nav4" slot).
# -------------------------------------------------------------------------

# Remember the old row 6 (nav4) English/French text before it gets
# overwritten - this text ("R code" / "Code R") moves down into the newly
# inserted row.
$oldEnglish = $ws.Range("B6").Value2
$oldFrench  = $ws.Range("C6").Value2

# Insert a new blank row directly beneath row 6 (i.e. above the old row 7),
# pushing every following row down by one.
$ws.Rows("7:7").Insert()

# Give the new row the normal data-row formatting (border/font) by copying
# it from the row immediately below, which still carries the original
# style used throughout the table.
$ws.Range("A8:C8").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 7: nav5 / R code / Code R
$ws.Range("A7").Value = "nav5"
$ws.Range("B7").Value = $oldEnglish
$ws.Range("C7").Value = $oldFrench

# Row 6 keeps the "nav4" id but now documents the new "Report" tab.
$ws.Range("B6").Value = "4. Report"
$ws.Range("C6").Value = "4. Rapport"

# The hyperlink-bearing "about" cells shifted from row 64 to row 65 along
# with everything else below the insertion point; re-anchor the two
# hyperlinks there (row insertion does not itself relocate them). Adding
# a hyperlink with a TextToDisplay would clobber the existing rich-text
# cell content, so add it address-only and then restore the original
# cell formatting (Hyperlinks.Add applies its own blue/underline style).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(65, 2), "https://doi.org/10.21105/joss.01082")
$ws.Hyperlinks.Add($ws.Cells.Item(65, 3), "https://doi.org/10.21105/joss.01082")
$ws.Range("B64:C64").Copy()
$ws.Range("B65:C65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final selection/scroll position shown in the diff.
$ws.Range("C6").Select()
